# "changes required for schnell added"
#
# The bulk-upload template's header row is being reworked: the old
# generic "Name"/"Number"/"Type" columns (for access cards) are dropped
# and replaced with dedicated vehicle columns, the card columns are
# moved after them, and a new "Device Name" column is appended.
#
# Resulting header row (row 1):
#   A1 Building Name   (unchanged)
#   B1 Floor Number    (unchanged)
#   C1 Flat Number     (unchanged)
#   D1 Vehicle Make    (new)
#   E1 Vehicle Number  (new)
#   F1 Vehicle Type    (new)
#   G1 Card Number     (moved from D1, text unchanged)
#   H1 Card Type       (moved from E1, text unchanged)
#   I1 Device Name     (new)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Vehicle Make"
$ws.Range("E1").Value = "Vehicle Number"
$ws.Range("F1").Value = "Vehicle Type"
$ws.Range("G1").Value = "Card Number"
$ws.Range("H1").Value = "Card Type"
$ws.Range("I1").Value = "Device Name"

# Match the cursor/selection left behind in the saved workbook.
$ws.Range("B3").Select() | Out-Null
